# Applies the LOT2004.docx edit: turns the two run-on, concatenated
# paragraphs ("Norma de recuperação" and "Bibliografia") into
# multi-line text using manual line breaks (w:br), matching how the
# sentences were evidently meant to be split.
$d = $word.ActiveDocument

# --- 1) "Norma de recuperação" paragraph -------------------------------
$old1 = "NF = (MF + PR)/2, onde PR é uma prova de recuperação.Prova de Recuperação (PR) para alunos com Média Final (MF) maior ou igual a 3,0 e menor doque 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final (NF) igual ou maiordo que 5,0."
$new1 = "NF = (MF + PR)/2, onde PR é uma prova de recuperação.^lProva de Recuperação (PR) para alunos com Média Final (MF) maior ou igual a 3,0 e menor do^lque 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final (NF) igual ou maior^ldo que 5,0."
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
if (-not $found1) {
    throw "Could not find the 'Norma de recuperacao' paragraph text to replace."
}
Write-Output "found1=$found1"

# --- 2) "Bibliografia" paragraph ----------------------------------------
$old2 = "1. Campbell, M.K.; Farrell, S.O. Bioquímica. Quinta edição. Editora Thomson Learning, São Paulo, 2008.2.Nelson, D.L.; Cox, M.M. Princípios de Bioquímica de Lehninger. Quinta Edition, Editora Artmed, Porto Alegre, 2011 3.Voet, D.; Voet, J.; Pratt, C.W. Bioquímica. Quarta Edição. Editora ARTMED, Porto Alegre, 2013 4.Berg, J.M., Tymoczko, J.L., Stryer, L. Bioquímica. Sétima edição. Editora Guanabara Koogan, Rio de Janeiro, 2014"
$new2 = "1. Campbell, M.K.; Farrell, S.O. Bioquímica. Quinta edição. Editora Thomson Learning, São Paulo, 2008.^l2.Nelson, D.L.; Cox, M.M. Princípios de Bioquímica de Lehninger. Quinta Edition, Editora Artmed, Porto Alegre, 2011 ^l3.Voet, D.; Voet, J.; Pratt, C.W. Bioquímica. Quarta Edição. Editora ARTMED, Porto Alegre, 2013 ^l4.Berg, J.M., Tymoczko, J.L., Stryer, L. Bioquímica. Sétima edição. Editora Guanabara Koogan, Rio de Janeiro, 2014"
$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
if (-not $found2) {
    throw "Could not find the 'Bibliografia' paragraph text to replace."
}
Write-Output "found2=$found2"
